# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect the latest scrape, per commit: "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): row -> new F value
$wsExhibit = $wb.Worksheets.Item("展览")
$exhibitUpdates = @{
    2  = 164
    3  = 41
    9  = 2256
    10 = 109
    13 = 1391
    14 = 492
    16 = 305
    19 = 24
    20 = 42
    23 = 61
    24 = 27
    25 = 1410
    27 = 362
    29 = 281
    30 = 352
}
foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Range("F$row").Value = $exhibitUpdates[$row]
}

# Sheet "全部类型" (sheet4): row -> new F value
$wsAll = $wb.Worksheets.Item("全部类型")
$allUpdates = @{
    2  = 164
    3  = 41
    10 = 2256
    11 = 109
    14 = 1391
    15 = 492
    17 = 305
    20 = 24
    21 = 42
    24 = 61
    25 = 27
    26 = 1410
    28 = 362
    30 = 281
    31 = 352
}
foreach ($row in $allUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allUpdates[$row]
}
